# Auto-generated edit script: update market-price derived cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1708.3572
$ws.Range("I40").Value = 1380
$ws.Range("J40").Value = 2146.1667
$ws.Range("K40").Value = 1380
$ws.Range("L40").Value = 2146.1667
$ws.Range("M40").Value = -1205
$ws.Range("N40").Value = -2496.1667
$ws.Range("H62").Value = 6317.4116
$ws.Range("I62").Value = 4700.4
$ws.Range("J62").Value = 8627.429
$ws.Range("K62").Value = 4700.4
$ws.Range("L62").Value = 8627.429
$ws.Range("M62").Value = -4076.4
$ws.Range("N62").Value = -9875.429
$ws.Range("H65").Value = 6317.4116
$ws.Range("I65").Value = 4700.4
$ws.Range("J65").Value = 8627.429
$ws.Range("K65").Value = 23502
$ws.Range("L65").Value = 43137.145
$ws.Range("M65").Value = -20382
$ws.Range("N65").Value = -49377.145
$ws.Range("H86").Value = 5906.095
$ws.Range("I86").Value = 1220.4286
$ws.Range("J86").Value = 15277.429
$ws.Range("K86").Value = 1220.4286
$ws.Range("L86").Value = 15277.429
$ws.Range("M86").Value = -97.42859999999996
$ws.Range("N86").Value = -17523.429
$ws.Range("H89").Value = 5906.095
$ws.Range("I89").Value = 1220.4286
$ws.Range("J89").Value = 15277.429
$ws.Range("K89").Value = 6102.143
$ws.Range("L89").Value = 76387.145
$ws.Range("M89").Value = -486.143
$ws.Range("N89").Value = -87619.145
$ws.Range("H116").Value = 3566.6875
$ws.Range("J116").Value = 4128.7
$ws.Range("L116").Value = 4128.7
$ws.Range("N116").Value = -11012.7
$ws.Range("H129").Value = 257394.64
$ws.Range("J129").Value = 323761.62
$ws.Range("L129").Value = 971284.86
$ws.Range("N129").Value = -981284.86
$ws.Range("H132").Value = 3143.7188
$ws.Range("I132").Value = 3584.1155
$ws.Range("K132").Value = 10752.3465
$ws.Range("M132").Value = -8222.3465
$ws.Range("I141").Value = 3331
$ws.Range("K141").Value = 9993
$ws.Range("M141").Value = -4813

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4624.3887
$ws.Range("I32").Value = 3952.5063
$ws.Range("J32").Value = 9449.727999999999
$ws.Range("K32").Value = 3952.5063
$ws.Range("L32").Value = 9449.727999999999
$ws.Range("M32").Value = -3665.5063
$ws.Range("N32").Value = -10023.728
$ws.Range("H114").Value = 31556.572
$ws.Range("J114").Value = 31556.572
$ws.Range("L114").Value = 31556.572
$ws.Range("N114").Value = -40234.572
$ws.Range("H119").Value = 30376.4
$ws.Range("J119").Value = 30376.4
$ws.Range("L119").Value = 30376.4
$ws.Range("N119").Value = -40052.4
$ws.Range("H124").Value = 7687.5
$ws.Range("J124").Value = 7687.5
$ws.Range("L124").Value = 7687.5
$ws.Range("N124").Value = -17507.5
$ws.Range("H125").Value = 34880
$ws.Range("J125").Value = 34880
$ws.Range("L125").Value = 34880
$ws.Range("N125").Value = -44720
$ws.Range("H135").Value = 37374.715
$ws.Range("J135").Value = 37374.715
$ws.Range("L135").Value = 37374.715
$ws.Range("N135").Value = -47514.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1112
$ws.Range("I94").Value = 918
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 918
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -467
$ws.Range("N94").Value = -2402
$ws.Range("H99").Value = 2249.75
$ws.Range("I99").Value = 2249.75
$ws.Range("K99").Value = 2249.75
$ws.Range("M99").Value = -751.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 40818.43
$ws.Range("J74").Value = 40818.43
$ws.Range("L74").Value = 40818.43
$ws.Range("N74").Value = -42566.43
$ws.Range("H77").Value = 40818.43
$ws.Range("J77").Value = 40818.43
$ws.Range("L77").Value = 122455.29
$ws.Range("N77").Value = -131191.29
$ws.Range("H88").Value = 18440
$ws.Range("J88").Value = 18440
$ws.Range("L88").Value = 18440
$ws.Range("N88").Value = -19252
$ws.Range("H91").Value = 18440
$ws.Range("J91").Value = 18440
$ws.Range("L91").Value = 18440
$ws.Range("N91").Value = -21248
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H134").Value = 1234.4166
$ws.Range("I134").Value = 1111.3
$ws.Range("J134").Value = 1850
$ws.Range("K134").Value = 3333.9
$ws.Range("L134").Value = 5550
$ws.Range("M134").Value = -798.8999999999996
$ws.Range("N134").Value = -10620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.05
$ws.Range("I2").Value = 18.294117
$ws.Range("K2").Value = 109.764702
$ws.Range("M2").Value = 3.235298
$ws.Range("H38").Value = 55555640
$ws.Range("I38").Value = 97.5
$ws.Range("K38").Value = 292.5
$ws.Range("M38").Value = 54.5
$ws.Range("H131").Value = 784.61957
$ws.Range("J131").Value = 779.1667
$ws.Range("L131").Value = 2337.5001
$ws.Range("N131").Value = -12417.5001
$ws.Range("H136").Value = 3296.2
$ws.Range("J136").Value = 4828.6665
$ws.Range("L136").Value = 14485.9995
$ws.Range("N136").Value = -24685.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 20005600
$ws.Range("J52").Value = 20005600
$ws.Range("L52").Value = 20005600
$ws.Range("N52").Value = -20006118
$ws.Range("H80").Value = 3149.4285
$ws.Range("I80").Value = 2753.5
$ws.Range("J80").Value = 3307.8
$ws.Range("K80").Value = 2753.5
$ws.Range("L80").Value = 3307.8
$ws.Range("M80").Value = -1755.5
$ws.Range("N80").Value = -5303.8
$ws.Range("H83").Value = 3149.4285
$ws.Range("I83").Value = 2753.5
$ws.Range("J83").Value = 3307.8
$ws.Range("K83").Value = 13767.5
$ws.Range("L83").Value = 16539
$ws.Range("M83").Value = -8775.5
$ws.Range("N83").Value = -26523
$ws.Range("H122").Value = 44445480
$ws.Range("I122").Value = 13889838
$ws.Range("J122").Value = 166668050
$ws.Range("K122").Value = 41669514
$ws.Range("L122").Value = 500004150
$ws.Range("M122").Value = -41667064
$ws.Range("N122").Value = -500009050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3197.742
$ws.Range("I40").Value = 2509.818
$ws.Range("J40").Value = 4879.3335
$ws.Range("K40").Value = 2509.818
$ws.Range("L40").Value = 4879.3335
$ws.Range("M40").Value = -2373.818
$ws.Range("N40").Value = -5151.3335
$ws.Range("H46").Value = 1241.421
$ws.Range("I46").Value = 1017.9375
$ws.Range("K46").Value = 1017.9375
$ws.Range("M46").Value = -829.9375
$ws.Range("H101").Value = 19000
$ws.Range("J101").Value = 19000
$ws.Range("L101").Value = 19000
$ws.Range("N101").Value = -25490
$ws.Range("H110").Value = 40006.5
$ws.Range("J110").Value = 40006.5
$ws.Range("L110").Value = 40006.5
$ws.Range("N110").Value = -48186.5
$ws.Range("H127").Value = 39858.707
$ws.Range("J127").Value = 39858.707
$ws.Range("L127").Value = 39858.707
$ws.Range("N127").Value = -49778.707

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 18250
$ws.Range("J69").Value = 18250
$ws.Range("L69").Value = 18250
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 18250
$ws.Range("J72").Value = 18250
$ws.Range("L72").Value = 54750
$ws.Range("N72").Value = -62238
$ws.Range("H82").Value = 15555
$ws.Range("J82").Value = 15555
$ws.Range("L82").Value = 15555
$ws.Range("N82").Value = -16321
$ws.Range("H85").Value = 15555
$ws.Range("J85").Value = 15555
$ws.Range("L85").Value = 15555
$ws.Range("N85").Value = -18207
